$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The weekly refresh re-shuffles the existing price records (rows 198-292,
# 95 rows x 18 cols) and appends 2 fresh ones, growing the sheet from
# A1:R292 to A1:R294:
#
#   new row 198         <- old row 215 (Choclero/Primera, Maule) with the
#                           date bumped to the new reporting date (44992)
#   new row 199         <- old row 216 (Choclero/Segunda, Maule) with the
#                           date bumped to the new reporting date (44992)
#   new rows 200..292    <- old rows 198..290  (shifted down by 2)
#   new row 293          <- old row 291
#   new row 294          <- old row 292
#
# Note: Range(...).Value2 returns a 1-based 2D array (COM semantics), while
# New-Object 'object[,]' returns a 0-based .NET array. Keep that straight
# below: $old is 1-based, $newData is 0-based.
# ---------------------------------------------------------------------------

$oldFirstRow = 198
$numCols     = 18  # columns A..R

$old = $ws.Range("A198:R292").Value2
$oldRowCount = $old.GetLength(0)

$newRowCount = $oldRowCount + 2
$newData = New-Object 'object[,]' $newRowCount, $numCols

for ($c = 1; $c -le $numCols; $c++) {
    # new row 198 (dest 0-based row 0) <- old row 215 (1-based src row 18)
    $newData[0, ($c - 1)] = $old[18, $c]
    # new row 199 (dest 0-based row 1) <- old row 216 (1-based src row 19)
    $newData[1, ($c - 1)] = $old[19, $c]
}
$newData[0, 3] = 44992
$newData[1, 3] = 44992

$lastShifted = $oldRowCount - 2
for ($i = 1; $i -le $lastShifted; $i++) {
    $dest = $i + 1
    for ($c = 1; $c -le $numCols; $c++) {
        $newData[$dest, ($c - 1)] = $old[$i, $c]
    }
}

# new row 293 (dest 0-based row 95) <- old row 291 (1-based src row 94)
# new row 294 (dest 0-based row 96) <- old row 292 (1-based src row 95)
for ($c = 1; $c -le $numCols; $c++) {
    $newData[95, ($c - 1)] = $old[94, $c]
    $newData[96, ($c - 1)] = $old[95, $c]
}

$newLastRow = $oldFirstRow + $newRowCount - 1
$targetRange = "A" + $oldFirstRow + ":R" + $newLastRow
$ws.Range($targetRange).Value2 = $newData

# The 2 brand-new rows (293, 294) sit beyond the sheet's previous used range,
# so their "Fecha" (date) cell needs the same date number-format the rest of
# column D already carries (style index 2 / numFmtId 165) - copy it from an
# existing date cell instead of leaving the new cells with the default format.
$dateFormat = $ws.Range("D197").NumberFormat
$ws.Range("D293").NumberFormat = $dateFormat
$ws.Range("D294").NumberFormat = $dateFormat
